# Release Form - F5: rename sheet to match the "5-Software Service Catalog -SS"
# form numbering (was copied from the Software Development Lifecycle forms).
#   F-SW-SD-05  ->  S-SW-SC-05
# The sheet-scoped Print_Area defined name embeds the sheet name in its
# RefersTo formula, so it needs to be refreshed too once the sheet is renamed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "S-SW-SC-05"

# Re-assert the print area (still A1:E13) so the Print_Area defined name
# re-serializes against the new sheet name instead of the stale one.
$ws.PageSetup.PrintArea = '$A$1:$E$13'
